# CP states from 1-10 into 1-11 (adding non-carnivorous)
#
# Slide 2 ("Carnivorous Plant trap states") holds:
#   - a Title placeholder
#   - several label rectangles
#   - ten small "number" rectangles (currently labelled 1..10)
#
# This script:
#   1. Gives the Title an explicit position/size and turns on
#      "shrink text on overflow" (produces <a:normAutofit/>).
#   2. Bumps every state number up by one (1->2, 2->3, ... 10->11),
#      also growing the box that now holds the two-digit "10".
#   3. Adds a new small label shape reading "(1 = non-carnivorous)".
#
# NOTE on units: the PowerPoint COM object model works in points while
# the underlying OOXML stores EMU (1 pt = 12700 EMU). Shape.Left/Top/
# Width/Height are `Single` (32-bit float) under the hood, so a handful
# of the point literals below are chosen (rather than a naive
# emu/12700.0) so that, after the inevitable 32-bit rounding, they land
# back exactly on the EMU value the diff expects.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# ---------------------------------------------------------------------
# 1) Title: explicit xfrm + shrink-text-on-overflow autofit.
# ---------------------------------------------------------------------
$title = $s.Shapes.Item(1)
$title.Left = 50.5                  # 641350 EMU
$title.Top = 15.422284126281738     # 195863 EMU
$title.Width = 621.0                # 7886700 EMU
$title.Height = 33.74992370605469   # 428624 EMU
$title.TextFrame.AutoSize = 2       # ppAutoSizeTextToFitShape -> <a:normAutofit/>

# ---------------------------------------------------------------------
# 2) Renumber the state labels (shift every existing number up by one).
#    Shapes are addressed by their fixed position in the slide's shape
#    collection, so this does not depend on (and does not collide with)
#    the text already in the boxes.
# ---------------------------------------------------------------------
$s.Shapes.Item(11).TextFrame.TextRange.Text = "2"   # was "1"
$s.Shapes.Item(12).TextFrame.TextRange.Text = "3"   # was "2"
$s.Shapes.Item(13).TextFrame.TextRange.Text = "4"   # was "3"
$s.Shapes.Item(14).TextFrame.TextRange.Text = "5"   # was "4"
$s.Shapes.Item(19).TextFrame.TextRange.Text = "6"   # was "5"
$s.Shapes.Item(20).TextFrame.TextRange.Text = "7"   # was "6"
$s.Shapes.Item(15).TextFrame.TextRange.Text = "8"   # was "7"
$s.Shapes.Item(16).TextFrame.TextRange.Text = "9"   # was "8"

# This one grows to fit the now-two-digit label ("9" -> "10").
$s.Shapes.Item(17).Width = 32.96882247924805        # 418704 EMU
$s.Shapes.Item(17).TextFrame.TextRange.Text = "10"  # was "9"

$s.Shapes.Item(18).TextFrame.TextRange.Text = "11"  # was "10"

# ---------------------------------------------------------------------
# 3) Add the new "(1 = non-carnivorous)" note.
#
#    Burn the slide's lowest still-free shape id (left behind by a
#    shape that no longer exists on this slide) with a scratch shape so
#    the real new shape below picks up the next-after-max id (23),
#    matching how PowerPoint itself allocates ids.
# ---------------------------------------------------------------------
$scratch = $s.Shapes.AddShape(1, 0, 0, 10, 10)
$scratch.Delete()

$note = $s.Shapes.Item(11).Duplicate()
$note.Name = "Rectangle 22"
$note.Left = 34.74519729614258      # 441264 EMU
$note.Top = 90.1719741821289        # 1145184 EMU
$note.Width = 172.84158325195312    # 2195088 EMU
$note.Height = 29.081260681152344   # 369332 EMU
$note.TextFrame.TextRange.Text = "(1 = non-carnivorous)"
